# Fix bug with item edit: the reservation for customer 10 (4-seat table,
# table count 3, dated 2024-04-29) was being duplicated instead of edited.
# Re-create both rows as they exist in the corrected workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Reservation ID 10
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(11, 3).Value = 45411
$ws.Cells.Item(11, 3).NumberFormat = "yyyy-MM-dd"
$ws.Cells.Item(11, 4).Value = 10
$ws.Cells.Item(11, 5).Value = "4 seat"
$ws.Cells.Item(11, 6).Value = 3

# Row 12 - Reservation ID 11 (duplicate entry caused by the edit bug)
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 10
$ws.Cells.Item(12, 3).Value = 45411
$ws.Cells.Item(12, 3).NumberFormat = "yyyy-MM-dd"
$ws.Cells.Item(12, 4).Value = 10
$ws.Cells.Item(12, 5).Value = "4 seat"
$ws.Cells.Item(12, 6).Value = 3
